# Updated the graph of execution times
# - refreshed Sheet1 timings (day 7 added, days 1-6 re-measured)
# - added Sheet2 with the raw per-day timing log lines
# - extended the 3 chart series (Initialization / Part 1 / Part 2) to cover
#   the new 7th day

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Refresh the existing day rows (2-7) on Sheet1 with the new timings
# ---------------------------------------------------------------------
$newData = @(
    @(1, 0.12790000000000001, 0.0191, 1.1194),
    @(2, 0.26569999999999999, 0.54059999999999997, 0),
    @(3, 1.4305000000000001, 0.0425, 0.16950000000000001),
    @(4, 0.43659999999999999, 12.691000000000001, 6.1646999999999998),
    @(5, 0.5242, 0, 0.0055),
    @(6, 2.4479000000000002, 1.2927, 0.59609999999999996),
    @(7, 0.14510000000000001, 0.0001, 0.0001)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = 2 + $i
    $vals = $newData[$i]
    $ws1.Range("A$row").Value = $vals[0]
    $ws1.Range("B$row").Value = $vals[1]
    $ws1.Range("C$row").Value = $vals[2]
    $ws1.Range("D$row").Value = $vals[3]
}

$ws1.Range("D2").Select()

# ---------------------------------------------------------------------
# 2. Add Sheet2 with the raw timing-log rows (one per day)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# columns I/J/K hold the same timings as Sheet1 B/C/D but in nanoseconds
$logRows = @(
    @(1, 127900, 19100, 1119400),
    @(2, 265700, 540600, 0),
    @(3, 1430500, 42500, 169500),
    @(4, 436600, 12691000, 6164700),
    @(5, 524200, 0, 5500),
    @(6, 2447900, 1292700, 596100),
    @(7, 145100, 100, 100)
)

for ($i = 0; $i -lt $logRows.Count; $i++) {
    $row = 1 + $i
    $vals = $logRows[$i]
    $ws2.Range("A$row").Value = "Day"
    $ws2.Range("B$row").Value = $vals[0]
    $ws2.Range("C$row").Value = "of"
    $ws2.Range("D$row").Value = 2020
    $ws2.Range("E$row").Value = "-"
    $ws2.Range("F$row").Value = "Timing"
    $ws2.Range("G$row").Value = "-"
    $ws2.Range("H$row").Value = "init:"
    $ws2.Range("I$row").Value = $vals[1]
    $ws2.Range("J$row").Value = $vals[2]
    $ws2.Range("K$row").Value = $vals[3]
    $ws2.Range("L$row").Value = "ns"
}

$ws2.Range("I1:K7").Select()

# ---------------------------------------------------------------------
# 3. Extend the 3 bar-chart series to the new 7-day range (A2:A8, etc.)
# ---------------------------------------------------------------------
$co = $ws1.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection()

$series.Item(1).Formula = "=SERIES(Sheet1!`$B`$1,Sheet1!`$A`$2:`$A`$8,Sheet1!`$B`$2:`$B`$8,1)"
$series.Item(2).Formula = "=SERIES(Sheet1!`$C`$1,Sheet1!`$A`$2:`$A`$8,Sheet1!`$C`$2:`$C`$8,2)"
$series.Item(3).Formula = "=SERIES(Sheet1!`$D`$1,Sheet1!`$A`$2:`$A`$8,Sheet1!`$D`$2:`$D`$8,3)"

$ws1.Select()
